$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove header-row (row 1) bold font / thin border / center+top alignment ---
# (the corresponding font/border/cellXf entries are dropped from styles.xml because
# no cell references them any more)
$headerRange = $ws.Range("A1:Q1")
$headerRange.ClearFormats()

# Helper: write a "blank but present" text cell (mirrors the workbook's pre-existing
# empty inline-string cells, e.g. B/H/L/P/Q) instead of simply clearing it, which
# would drop the cell from the sheet entirely.
function Set-BlankCell($addr) {
    $ws.Range($addr).Value = "'"
    $ws.Range($addr).ClearFormats()
}

# --- A1: was "Unnamed: 0", now blank ---
Set-BlankCell "A1"

# --- Corrected pre/post/total fixation metrics (rows 3-7) ---
# Columns I, M, O are cleared out (no longer populated) for every metric row;
# C, D, F, G, J, K, N get corrected values; E is left untouched.

$newValues = @{
    "C3" = 15;      "D3" = 1;        "F3" = 2;       "G3" = 2;       "J3" = 4;        "K3" = 1;       "N3" = 67
    "C4" = 28;      "D4" = 2;        "F4" = 3;       "G4" = 3;       "J4" = 5;        "K4" = 2;       "N4" = 363
    "C5" = 6306.46; "D5" = 283.52;   "F5" = 608.78;  "G5" = 583.83;  "J5" = 1201.23;  "K5" = 667.36;  "N5" = 116082.3
    "C6" = 2.02;    "D6" = 0.09;     "F6" = 0.2;     "G6" = 0.19;    "J6" = 0.38;     "K6" = 0.21;    "N6" = 37.19
    "C7" = 225.23;  "D7" = 141.76;   "F7" = 202.93;  "G7" = 194.61;  "J7" = 240.25;   "K7" = 333.68;  "N7" = 319.79
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}

foreach ($row in 3..7) {
    foreach ($col in @("I", "M", "O")) {
        Set-BlankCell "$col$row"
    }
}
